$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "searchTable(tableNo)" and "getTable(tableNo)" operation-contract
# blocks (previously rows 25-40) entirely; rows below shift up to close the gap.
$ws.Range("A25:A40").EntireRow.Delete() | Out-Null

# Remove the extra column H mini-table that sat next to the first block.
$ws.Range("H1:H10").EntireColumn.Delete() | Out-Null

# ---- Block 1 (makeReservation) now at rows 1-7: fill in the two new
# post-condition lines that used to be blank.
$ws.Range("B5").Value = "New reservation instance is created"
$ws.Range("B6").Value = "New reservation is saved in the database"

# ---- Block 2 (insertTableNo -> enterTableNo) now at rows 9-15
$ws.Range("B9").Value = "enterTableNo(tableNo)"
$ws.Range("B12").Value = "table number was entered"

# ---- Block 3 (findTableByNo) now at rows 17-23
$ws.Range("B19").Value = "table number was entered"
$ws.Range("B20").Value = "The system finds the table"

# Update the selected cell (the sheet view also no longer needs to be
# scrolled down to row 49, since the content now fits in rows 1-55).
$ws.Range("F8").Select() | Out-Null
